$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1) "Three fingers up " + "– Moves Right" (two runs) -> merge into a
#    single run "Three fingers up – Moves Right" that keeps the
#    "dirty" run's formatting (the second run's rPr survives).
# ---------------------------------------------------------------------
$gestures = $s.Shapes.Item("Content Placeholder 1")
$gTextRange = $gestures.TextFrame.TextRange

$gFull = $gTextRange.Text
$firstPhrase = "Three fingers up "
$secondPhrase = "– Moves Right"
$gIdx = $gFull.IndexOf($firstPhrase + $secondPhrase)

# Clear the first run's text (it disappears, leaving only the second run)
$firstRun = $gTextRange.Characters($gIdx + 1, $firstPhrase.Length)
$firstRun.Text = ""

# Re-find the remaining run and prepend the full combined text onto it,
# so the merged text ends up carried by what used to be the second run
# (and therefore keeps its rPr, e.g. dirty="0").
$gFull2 = $gTextRange.Text
$gIdx2 = $gFull2.IndexOf($secondPhrase)
$secondRun = $gTextRange.Characters($gIdx2 + 1, $secondPhrase.Length)
$secondRun.Text = $firstPhrase + $secondPhrase

# ---------------------------------------------------------------------
# 2) ", Archie Wills 47440188, Liam Hugo 4749048" -> split into three
#    runs: ", Archie Wills 47440188, " / "Liam Ryan " / "4749048"
#    ("Liam Hugo" doesn't exist -> "Liam Ryan").
# ---------------------------------------------------------------------
$subtitle = $s.Shapes.Item("Subtitle 3")
$sTextRange = $subtitle.TextFrame.TextRange

$sFull = $sTextRange.Text
$needle = ", Archie Wills 47440188, Liam Hugo 4749048"
$base = $sFull.IndexOf($needle)

$seg1 = ", Archie Wills 47440188, "
$seg2old = "Liam Hugo "
$seg2new = "Liam Ryan "
$seg3 = "4749048"

# Work right-to-left so earlier offsets stay valid.
$seg3Start = $base + $seg1.Length + $seg2old.Length
$seg3Range = $sTextRange.Characters($seg3Start + 1, $seg3.Length)
$seg3Range.Text = $seg3

$seg2Start = $base + $seg1.Length
$seg2Range = $sTextRange.Characters($seg2Start + 1, $seg2old.Length)
$seg2Range.Text = $seg2new

$seg1Start = $base
$seg1Range = $sTextRange.Characters($seg1Start + 1, $seg1.Length)
$seg1Range.Text = $seg1
